# Update the "fromDate" and "toDate" value cells on the worksheet.
# These cells are stored as text (quote-prefixed numbers), so we use a
# leading apostrophe to force a text value and keep the existing style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "'27"
$ws.Range("B8").Value = "'29"

# Update the active selection to C9, matching the saved view state.
$ws.Range("C9").Select()
